$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 390 (pushes existing rows 390.. down to 392..)
$ws.Rows(390).Insert()
$ws.Rows(390).Insert()

# Row 390 - "Primera" record for the new date 44522
$ws.Cells.Item(390, 1).Value = 8
$ws.Cells.Item(390, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(390, 3).Value = "Coquimbo"
$ws.Cells.Item(390, 4).Value = 44522
$ws.Cells.Item(390, 5).Value = 4
$ws.Cells.Item(390, 6).Value = 100112008
$ws.Cells.Item(390, 7).Value = "Coliflor"
$ws.Cells.Item(390, 8).Value = "Sin especificar"
$ws.Cells.Item(390, 9).Value = "Primera"
$ws.Cells.Item(390, 10).Value = 2600
$ws.Cells.Item(390, 11).Value = 600
$ws.Cells.Item(390, 12).Value = 700
$ws.Cells.Item(390, 13).Value = 650
$ws.Cells.Item(390, 14).Value = "$/unidad"
$ws.Cells.Item(390, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(390, 16).Value = 650
$ws.Cells.Item(390, 17).Value = 1
$ws.Cells.Item(390, 18).Value = "Hortaliza"

# Row 391 - "Segunda" record for the new date 44522
$ws.Cells.Item(391, 1).Value = 8
$ws.Cells.Item(391, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(391, 3).Value = "Coquimbo"
$ws.Cells.Item(391, 4).Value = 44522
$ws.Cells.Item(391, 5).Value = 4
$ws.Cells.Item(391, 6).Value = 100112008
$ws.Cells.Item(391, 7).Value = "Coliflor"
$ws.Cells.Item(391, 8).Value = "Sin especificar"
$ws.Cells.Item(391, 9).Value = "Segunda"
$ws.Cells.Item(391, 10).Value = 1600
$ws.Cells.Item(391, 11).Value = 500
$ws.Cells.Item(391, 12).Value = 550
$ws.Cells.Item(391, 13).Value = 525
$ws.Cells.Item(391, 14).Value = "$/unidad"
$ws.Cells.Item(391, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(391, 16).Value = 525
$ws.Cells.Item(391, 17).Value = 1
$ws.Cells.Item(391, 18).Value = "Hortaliza"

# Fix a pre-existing rounding value (674 -> 675) on the row that is now 443
# (was row 441 before the two rows above were inserted)
$ws.Cells.Item(443, 13).Value = 675
$ws.Cells.Item(443, 16).Value = 675
